$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# A new weekly price record was added at the top of the "Haba" data block
# (row 100). This pushes all existing records (old rows 100-135) down by
# one row (new rows 101-136), preserving the existing last row's data in
# the newly created row 136.
$ws.Rows("100:100").Insert()

# Populate the newly inserted row 100 with the new record's data.
$ws.Cells.Item(100, 1).Value = 4
$ws.Cells.Item(100, 2).Value = "Feria Lagunitas de Puerto Montt"
$ws.Cells.Item(100, 3).Value = "Los Lagos"
$ws.Cells.Item(100, 4).Value = 45093
$ws.Cells.Item(100, 5).Value = 10
$ws.Cells.Item(100, 6).Value = 100112026
$ws.Cells.Item(100, 7).Value = "Haba"
$ws.Cells.Item(100, 8).Value = "Sin especificar"
$ws.Cells.Item(100, 9).Value = "Primera"
$ws.Cells.Item(100, 10).Value = 80
$ws.Cells.Item(100, 11).Value = 21000
$ws.Cells.Item(100, 12).Value = 21000
$ws.Cells.Item(100, 13).Value = 21000
$ws.Cells.Item(100, 14).Value = "`$/saco 25 kilos"
$ws.Cells.Item(100, 15).Value = "Provincia de Limarí"
$ws.Cells.Item(100, 16).Value = 840
$ws.Cells.Item(100, 17).Value = 25
$ws.Cells.Item(100, 18).Value = "Hortaliza"
